$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: restructure text and center it.
#    Old: "Safeguarding and Child Safeguarding Questions incorporated on the
#          Applicant Tracking System (ATS)"  (with a comment anchored on "ATS")
#    New: "Safeguarding (Child and Adult Safeguarding) Questions incorporated
#          on the Applicant Tracking System (ATS)"
# ---------------------------------------------------------------------------

# Remove the existing comment anchored to "ATS" in the title.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments($i).Delete()
}

$p1 = $d.Paragraphs(1).Range
$p1Start = $p1.Start
$p1End = $p1.End
$titleRange = $d.Range($p1Start, $p1End - 1)
$titleRange.Text = "Safeguarding (Child and Adult Safeguarding) Questions incorporated on the Applicant Tracking System (ATS)"
$d.Paragraphs(1).Range.ParagraphFormat.Alignment = 1

# ---------------------------------------------------------------------------
# 2. "Have you ever been charged with a criminal offence" question ->
#    "Have you ever been charged with a sexual exploitation or sexual
#     harassment offence"
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Have you ever been charged with a criminal offence", $true, $false, $false, $false, $false, $true, 1, $false, "Have you ever been charged with a sexual exploitation or sexual harassment offence", 2)

# ---------------------------------------------------------------------------
# 3. "Have you ever received a written warning ... unsatisfactory work
#     performance?" -> "Have you ever received a written or verbal warning
#     ... in relation to safeguarding violations?"
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Have you ever received a written warning or been dismissed or resigned following allegations of improper or unprofessional conduct or unsatisfactory work performance?", $true, $false, $false, $false, $false, $true, 1, $false, "Have you ever received a written or verbal warning or been dismissed or resigned following allegations of improper or unprofessional conduct in relation to safeguarding violations?", 2)

# ---------------------------------------------------------------------------
# 4. Referees paragraph: reword and split off the "If yes, please provide
#    details" sentence into its own paragraph (indented, justified).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Safeguarding, Child Safeguarding, Prevention from Sexual Exploitation and Abuse ", $true, $false, $false, $false, $false, $true, 1, $false, "Adult and Child Safeguarding, Sexual Exploitation and Abuse ", 2)

$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Do you foresee any problem arising from this process? If yes, please provide details: …………. ", $true, $false, $false, $false, $false, $true, 1, $false, "Do you foresee any problem arising from this process? ^pIf yes, please provide details: …………. ", 2)

# Remove the now-orphaned _GoBack bookmark left inside the referees paragraph
# (it will be re-created a little further down, right before the final
# paragraphs).
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# Locate the referees paragraph and the freshly split-off "If yes" paragraph
# so we can fix up paragraph formatting to match the target.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if ($ptext -like "If yes, please provide details: …………. *") {
        $ifYesIndex = $i
    }
}

$ifYesPara = $d.Paragraphs($ifYesIndex)
$ifYesPara.Range.ParagraphFormat.LeftIndent = 18
$ifYesPara.Range.ParagraphFormat.Alignment = 3

# ---------------------------------------------------------------------------
# 5. Insert a new bookmark paragraph (_GoBack) followed by one more blank
#    paragraph at the end of the document.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$endOfDoc = $d.Paragraphs($lastParaIndex).Range
$endOfDoc.InsertParagraphAfter()

$lastParaIndex = $d.Paragraphs.Count
$d.Bookmarks.Add("_GoBack", $d.Paragraphs($lastParaIndex - 1).Range)

$endOfDoc = $d.Paragraphs($d.Paragraphs.Count).Range
$endOfDoc.InsertParagraphAfter()
